# Scheduled-runner update: refresh cached market-board prices / leve profit figures
# across the Asura_Profits sheets (one block of edits per affected leve row).
$wb = $excel.ActiveWorkbook

# ALC row 15 (Leve Item ID 44146)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 1500.3594
$ws.Range("I15").Value = 1500.3594
$ws.Range("K15").Value = 4501.0782
$ws.Range("M15").Value = -4332.0782

# ALC row 98 (Leve Item ID 36237)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H98").Value = 4714.1724
$ws.Range("I98").Value = 3102.5
$ws.Range("J98").Value = 9779.429
$ws.Range("K98").Value = 3102.5
$ws.Range("L98").Value = 9779.429
$ws.Range("M98").Value = -1604.5
$ws.Range("N98").Value = -12775.429

# ALC row 100 (Leve Item ID 19906)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H100").Value = 3147.647
$ws.Range("I100").Value = 2617.5
$ws.Range("J100").Value = 3436.818
$ws.Range("K100").Value = 2617.5
$ws.Range("L100").Value = 3436.818
$ws.Range("M100").Value = -2076.5
$ws.Range("N100").Value = -4518.818

# ALC row 122 (Leve Item ID 36237)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H122").Value = 4714.1724
$ws.Range("I122").Value = 3102.5
$ws.Range("J122").Value = 9779.429
$ws.Range("K122").Value = 9307.5
$ws.Range("L122").Value = 29338.287
$ws.Range("M122").Value = -6857.5
$ws.Range("N122").Value = -34238.287

# ALC row 125 (Leve Item ID 36228)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H125").Value = 71901.42999999999
$ws.Range("J125").Value = 512.1
$ws.Range("L125").Value = 4608.900000000001
$ws.Range("N125").Value = -9528.900000000001

# ALC row 127 (Leve Item ID 36114)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H127").Value = 968.41
$ws.Range("I127").Value = 398.6
$ws.Range("J127").Value = 998.4
$ws.Range("K127").Value = 1195.8
$ws.Range("L127").Value = 2995.2
$ws.Range("M127").Value = 3764.2
$ws.Range("N127").Value = -12915.2

# ARM row 45 (Leve Item ID 27714)
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 1239.3077
$ws.Range("I45").Value = 1175.9166
$ws.Range("J45").Value = 2000
$ws.Range("K45").Value = 1175.9166
$ws.Range("L45").Value = 2000
$ws.Range("M45").Value = -798.9166
$ws.Range("N45").Value = -2754

# ARM row 102 (Leve Item ID 19945)
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H102").Value = 0
$ws.Range("I102").Value = 0
$ws.Range("J102").Value = 0
$ws.Range("K102").Value = 0
$ws.Range("L102").Value = ""
$ws.Range("M102").Value = ""
$ws.Range("N102").Value = 0

# ARM row 122 (Leve Item ID 36168)
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H122").Value = 1908.381
$ws.Range("I122").Value = 1775.8572
$ws.Range("J122").Value = 2173.4285
$ws.Range("K122").Value = 5327.571599999999
$ws.Range("L122").Value = 6520.2855
$ws.Range("M122").Value = -2877.571599999999
$ws.Range("N122").Value = -11420.2855

# BSM row 105 (Leve Item ID 19947)
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 3349.1667
$ws.Range("I105").Value = 3073
$ws.Range("J105").Value = 3901.5
$ws.Range("K105").Value = 3073
$ws.Range("L105").Value = 3901.5
$ws.Range("M105").Value = -1326
$ws.Range("N105").Value = -7395.5

# CRP row 99 (Leve Item ID 36198)
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H99").Value = 2779.9
$ws.Range("I99").Value = 2778.9333
$ws.Range("J99").Value = 2782.8
$ws.Range("K99").Value = 2778.9333
$ws.Range("L99").Value = 2782.8
$ws.Range("M99").Value = -1280.9333
$ws.Range("N99").Value = -5778.8

# CRP row 126 (Leve Item ID 36198)
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H126").Value = 2779.9
$ws.Range("I126").Value = 2778.9333
$ws.Range("J126").Value = 2782.8
$ws.Range("K126").Value = 8336.7999
$ws.Range("L126").Value = 8348.400000000001
$ws.Range("M126").Value = -5866.7999
$ws.Range("N126").Value = -13288.4

# CRP row 138 (Leve Item ID 42302)
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H138").Value = 0
$ws.Range("J138").Value = 0
$ws.Range("L138").Value = ""
$ws.Range("N138").Value = 0

# CUL row 5 (Leve Item ID 43974)
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 2119.9333
$ws.Range("I5").Value = 2672.111
$ws.Range("K5").Value = 8016.333
$ws.Range("M5").Value = -7904.333

# CUL row 34 (Leve Item ID 4749)
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H34").Value = 608.7143
$ws.Range("J34").Value = 640
$ws.Range("L34").Value = 1920
$ws.Range("N34").Value = -2088

# CUL row 39 (Leve Item ID 4712)
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H39").Value = 2416
$ws.Range("J39").Value = 2416
$ws.Range("L39").Value = 7248
$ws.Range("N39").Value = -7836

# CUL row 55 (Leve Item ID 4733)
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H55").Value = 4322.4
$ws.Range("J55").Value = 4481.6665
$ws.Range("L55").Value = 13444.9995
$ws.Range("N55").Value = -13798.9995

# CUL row 107 (Leve Item ID 27838)
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H107").Value = 461.33334
$ws.Range("I107").Value = 539.55554
$ws.Range("J107").Value = 226.66667
$ws.Range("K107").Value = 1618.66662
$ws.Range("L107").Value = 680.00001
$ws.Range("M107").Value = 301.33338
$ws.Range("N107").Value = -4520.00001

# CUL row 135 (Leve Item ID 43974)
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H135").Value = 2119.9333
$ws.Range("I135").Value = 2672.111
$ws.Range("K135").Value = 24048.999
$ws.Range("M135").Value = -21513.999

# GSM row 97 (Leve Item ID 19940)
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 56876.668
$ws.Range("I97").Value = 56876.668
$ws.Range("K97").Value = 56876.668
$ws.Range("M97").Value = -56380.668

# GSM row 122 (Leve Item ID 36182)
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 2674
$ws.Range("I122").Value = 1162.8
$ws.Range("J122").Value = 3933.3333
$ws.Range("K122").Value = 3488.4
$ws.Range("L122").Value = 11799.9999
$ws.Range("M122").Value = -1038.4
$ws.Range("N122").Value = -16699.9999

# GSM row 123 (Leve Item ID 34150)
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H123").Value = 13828.765
$ws.Range("J123").Value = 13828.765
$ws.Range("L123").Value = 13828.765
$ws.Range("N123").Value = -18728.765

# LTW row 7 (Leve Item ID 36249)
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2943.8667
$ws.Range("I7").Value = 2844.125
$ws.Range("K7").Value = 2844.125
$ws.Range("M7").Value = -2732.125

# LTW row 40 (Leve Item ID 36248)
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 3592.8667
$ws.Range("I40").Value = 3441.1428
$ws.Range("K40").Value = 3441.1428
$ws.Range("M40").Value = -3305.1428

# LTW row 93 (Leve Item ID 19993)
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H93").Value = 2000
$ws.Range("I93").Value = 0
$ws.Range("J93").Value = 2000
$ws.Range("K93").Value = 0
$ws.Range("L93").Value = ""
$ws.Range("M93").Value = 2000
$ws.Range("N93").Value = -4496

# LTW row 122 (Leve Item ID 36247)
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 23689190
$ws.Range("I122").Value = 17861614
$ws.Range("K122").Value = 53584842
$ws.Range("M122").Value = -53582392

# LTW row 126 (Leve Item ID 36249)
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H126").Value = 2943.8667
$ws.Range("I126").Value = 2844.125
$ws.Range("K126").Value = 8532.375
$ws.Range("M126").Value = -6062.375

# LTW row 139 (Leve Item ID 43310)
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H139").Value = 52305
$ws.Range("J139").Value = 58457.5
$ws.Range("L139").Value = 58457.5
$ws.Range("N139").Value = -68737.5

# WVR row 21 (Leve Item ID 3341)
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H21").Value = 50000
$ws.Range("J21").Value = 0
$ws.Range("L21").Value = 0
$ws.Range("N21").Value = ""

# WVR row 35 (Leve Item ID 3341)
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H35").Value = 50000
$ws.Range("J35").Value = 0
$ws.Range("L35").Value = 0
$ws.Range("N35").Value = ""

# WVR row 122 (Leve Item ID 36208)
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 10778263
$ws.Range("I122").Value = 11906881
$ws.Range("J122").Value = 7815640
$ws.Range("K122").Value = 35720643
$ws.Range("L122").Value = 23446920
$ws.Range("M122").Value = -35718193
$ws.Range("N122").Value = -23451820

# WVR row 123 (Leve Item ID 34127)
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H123").Value = 21861
$ws.Range("J123").Value = 21861
$ws.Range("L123").Value = 21861
$ws.Range("N123").Value = -31661

# WVR row 126 (Leve Item ID 36210)
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 10357.526
$ws.Range("I126").Value = 10357.526
$ws.Range("K126").Value = 31072.578
$ws.Range("M126").Value = -28602.578

# WVR row 136 (Leve Item ID 44031)
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 1620.3914
$ws.Range("I136").Value = 1635.2632
$ws.Range("J136").Value = 1549.75
$ws.Range("K136").Value = 4905.7896
$ws.Range("L136").Value = 4649.25
$ws.Range("M136").Value = -2355.7896
$ws.Range("N136").Value = -9749.25
